$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepare number formats for the new row range so that Excel's automatic
# type-detection doesn't turn numeric-looking / date-looking text into
# numbers or date serials, and so that the '#' / 'Contact' columns keep
# the same integer number format used by the existing data rows (style s="65").
$ws.Range("A57:A71").NumberFormat = "0"
$ws.Range("D57:D71").NumberFormat = "0"
$ws.Range("B57:B71").NumberFormat = "@"
$ws.Range("E57:E71").NumberFormat = "@"

# Row 57 (# 55)
$ws.Range("A57").Value = 55
$ws.Range("B57").Value = '22-12-2025'
$ws.Range("C57").Value = 'ajith'
$ws.Range("D57").Value = 8943409950
$ws.Range("E57").Value = '04-04-2026'
$ws.Range("F57").Value = 'SHAHEEN'
$ws.Range("G57").Value = 'Loss'
$ws.Range("H57").Value = 'ENQUIRY'
$ws.Range("I57").Value = 'ENQUIRY WITHOUT TRIAL'
$ws.Range("J57").Value = '-'
$ws.Range("K57").Value = 'CHECKING FOR PRICE AND COLLECTION'

# Row 58 (# 56)
$ws.Range("A58").Value = 56
$ws.Range("B58").Value = '22-12-2025'
$ws.Range("C58").Value = 'aslam'
$ws.Range("D58").Value = 7994224600
$ws.Range("E58").Value = '28-12-2025'
$ws.Range("F58").Value = 'MUHAMMED RAFI P V'
$ws.Range("G58").Value = 'Loss'
$ws.Range("H58").Value = 'ENQUIRY'
$ws.Range("I58").Value = 'ENQUIRY WITHOUT TRIAL'
$ws.Range("J58").Value = '-'
$ws.Range("K58").Value = 'collections and prices'

# Row 59 (# 57)
$ws.Range("A59").Value = 57
$ws.Range("B59").Value = '22-12-2025'
$ws.Range("C59").Value = 'anshif'
$ws.Range("D59").Value = 7736155435
$ws.Range("E59").Value = '29-12-2025'
$ws.Range("F59").Value = 'AKSHAY. V'
$ws.Range("G59").Value = 'Loss'
$ws.Range("H59").Value = 'CUSTOMER INTERNAL ISSUES'
$ws.Range("I59").Value = 'FAMILY DISAPPROVEL'
$ws.Range("J59").Value = '-'
$ws.Range("K59").Value = 'CUSTOMER WILL DECIDE TMRW'

# Row 60 (# 58)
$ws.Range("A60").Value = 58
$ws.Range("B60").Value = '22-12-2025'
$ws.Range("C60").Value = 'MUSTHAQ'
$ws.Range("D60").Value = 7994165151
$ws.Range("E60").Value = '03-01-2026'
$ws.Range("F60").Value = 'MUHAMMED ANSHIF C.K'
$ws.Range("G60").Value = 'Loss'
$ws.Range("H60").Value = 'CUSTOMER INTERNAL ISSUES'
$ws.Range("I60").Value = 'FAMILY DISAPPROVEL'
$ws.Range("J60").Value = '-'
$ws.Range("K60").Value = 'WILL DECIDE TMRW'

# Row 61 (# 59)
$ws.Range("A61").Value = 59
$ws.Range("B61").Value = '24-12-2025'
$ws.Range("C61").Value = 'vishnu'
$ws.Range("D61").Value = 6235500369
$ws.Range("E61").Value = '28-12-2025'
$ws.Range("F61").Value = 'SREEJESH C S'
$ws.Range("G61").Value = 'Loss'
$ws.Range("H61").Value = 'ENQUIRY'
$ws.Range("I61").Value = 'ENQUIRY WITHOUT TRIAL'
$ws.Range("J61").Value = '-'
$ws.Range("K61").Value = 'loss'

# Row 62 (# 60)
$ws.Range("A62").Value = 60
$ws.Range("B62").Value = '24-12-2025'
$ws.Range("C62").Value = 'FARHAN'
$ws.Range("D62").Value = 7034108801
$ws.Range("E62").Value = '27-12-2025'
$ws.Range("F62").Value = 'MUHAMMED ANSHIF C.K'
$ws.Range("G62").Value = 'Loss'
$ws.Range("H62").Value = 'PRODUCT'
$ws.Range("I62").Value = 'REQUIRED MODEL NOT AVAILABLE'
$ws.Range("J62").Value = '-'
$ws.Range("K62").Value = 'CUSTOMER NEEDED BLACK SUIT WITH HEAVY STONE WORK'

# Row 63 (# 61)
$ws.Range("A63").Value = 61
$ws.Range("B63").Value = '24-12-2025'
$ws.Range("C63").Value = 'amjith'
$ws.Range("D63").Value = 8086437713
$ws.Range("E63").Value = '08-01-2026'
$ws.Range("F63").Value = 'MUHAMMED RAFI P V'
$ws.Range("G63").Value = 'Loss'
$ws.Range("H63").Value = 'CUSTOMER INTERNAL ISSUES'
$ws.Range("I63").Value = 'FAMILY DISAPPROVEL'
$ws.Range("J63").Value = '-'
$ws.Range("K63").Value = 'they will decide after confirming with their family in two days'

# Row 64 (# 62)
$ws.Range("A64").Value = 62
$ws.Range("B64").Value = '24-12-2025'
$ws.Range("C64").Value = 'asarudheen'
$ws.Range("D64").Value = 7025525525
$ws.Range("E64").Value = '27-12-2025'
$ws.Range("F64").Value = 'MUHAMMED RAFI P V'
$ws.Range("G64").Value = 'Loss'
$ws.Range("H64").Value = 'CUSTOMER INTERNAL ISSUES'
$ws.Range("I64").Value = 'FAMILY DISAPPROVEL'
$ws.Range("J64").Value = '-'
$ws.Range("K64").Value = 'tomorrow will confirm'

# Row 65 (# 63)
$ws.Range("A65").Value = 63
$ws.Range("B65").Value = '24-12-2025'
$ws.Range("C65").Value = 'ashif'
$ws.Range("D65").Value = 9061234273
$ws.Range("E65").Value = '04-01-2026'
$ws.Range("F65").Value = 'SHAHEEN'
$ws.Range("G65").Value = 'Loss'
$ws.Range("H65").Value = 'PRODUCT'
$ws.Range("I65").Value = 'REQUIRED MODEL NOT AVAILABLE'
$ws.Range("J65").Value = '-'
$ws.Range("K65").Value = 'full work'

# Row 66 (# 64)
$ws.Range("A66").Value = 64
$ws.Range("B66").Value = '25-12-2025'
$ws.Range("C66").Value = 'iqbal'
$ws.Range("D66").Value = 8589378337
$ws.Range("E66").Value = '28-12-2025'
$ws.Range("F66").Value = 'SHAHEEN'
$ws.Range("G66").Value = 'Loss'
$ws.Range("H66").Value = 'ENQUIRY'
$ws.Range("I66").Value = 'ENQUIRY WITHOUT TRIAL'
$ws.Range("J66").Value = '-'
$ws.Range("K66").Value = 'CHECKING FOR PRICE AND COLLECTION'

# Row 67 (# 65)
$ws.Range("A67").Value = 65
$ws.Range("B67").Value = '25-12-2025'
$ws.Range("C67").Value = 'Adhil'
$ws.Range("D67").Value = 9747291612
$ws.Range("E67").Value = '27-12-2025'
$ws.Range("F67").Value = 'AKSHAY. V'
$ws.Range("G67").Value = 'Loss'
$ws.Range("H67").Value = 'ENQUIRY'
$ws.Range("I67").Value = 'Enquiry for Relative/Friend'
$ws.Range("J67").Value = '-'
$ws.Range("K67").Value = 'double breasted'

# Row 68 (# 66)
$ws.Range("A68").Value = 66
$ws.Range("B68").Value = '25-12-2025'
$ws.Range("C68").Value = 'MANAS'
$ws.Range("D68").Value = 9447424128
$ws.Range("E68").Value = '17-01-2026'
$ws.Range("F68").Value = 'MUHAMMED ANSHIF C.K'
$ws.Range("G68").Value = 'Loss'
$ws.Range("H68").Value = 'CUSTOMER INTERNAL ISSUES'
$ws.Range("I68").Value = 'FAMILY DISAPPROVEL'
$ws.Range("J68").Value = '-'
$ws.Range("K68").Value = 'WILL FINALIZE TOMORROW'

# Row 69 (# 67)
$ws.Range("A69").Value = 67
$ws.Range("B69").Value = '25-12-2025'
$ws.Range("C69").Value = 'sadiq'
$ws.Range("D69").Value = 9744566707
$ws.Range("E69").Value = '05-01-2026'
$ws.Range("F69").Value = 'MUHAMMED RAFI P V'
$ws.Range("G69").Value = 'Loss'
$ws.Range("H69").Value = 'CUSTOMER INTERNAL ISSUES'
$ws.Range("I69").Value = 'FAMILY DISAPPROVEL'
$ws.Range("J69").Value = '-'
$ws.Range("K69").Value = 'will decide after two days confirm'

# Row 70 (# 68)
$ws.Range("A70").Value = 68
$ws.Range("B70").Value = '25-12-2025'
$ws.Range("C70").Value = 'shabu'
$ws.Range("D70").Value = 9526561242
$ws.Range("E70").Value = '27-12-2025'
$ws.Range("F70").Value = 'MUHAMMED RAFI P V'
$ws.Range("G70").Value = 'Loss'
$ws.Range("H70").Value = 'SIZE NOT SUITABLE'
$ws.Range("I70").Value = 'SIZE TOO SMALL'
$ws.Range("J70").Value = '-'
$ws.Range("K70").Value = 'customer needed indowestern model but size too small'

# Row 71 (# 69)
$ws.Range("A71").Value = 69
$ws.Range("B71").Value = '25-12-2025'
$ws.Range("C71").Value = 'nihal'
$ws.Range("D71").Value = 8086202129
$ws.Range("E71").Value = '19-01-2026'
$ws.Range("F71").Value = 'MUHAMMED RAFI P V'
$ws.Range("G71").Value = 'Loss'
$ws.Range("H71").Value = 'CUSTOMER INTERNAL ISSUES'
$ws.Range("I71").Value = 'FAMILY DISAPPROVEL'
$ws.Range("J71").Value = '-'
$ws.Range("K71").Value = 'they will decide after confirming with their family in two days'
